$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("enemy")
$ws2 = $wb.Worksheets.Item("enemy_group")

# ---------------------------------------------------------------
# Sheet "enemy": tweak a handful of dice/parameter columns (F:K)
# for the BAT / BAT_2 / WOLF rows (rows 3-6).
# ---------------------------------------------------------------

# Row 3 (BAT)
$ws1.Range("H3").Value = 2

# Row 4 (BAT_2)
$ws1.Range("H4").Value = 3
$ws1.Range("I4").Value = 7

# Row 5 (BAT_3) - clear all of F5:K5
$ws1.Range("F5:K5").ClearContents()

# Row 6 (WOLF)
$ws1.Range("F6").Value = 10
$ws1.Range("G6").Value = 10
$ws1.Range("H6").Value = 2

# ---------------------------------------------------------------
# Sheet "enemy_group": re-point some group members and blank out
# the slots that no longer hold a member. The emptied cells get
# switched to the plain Verdana font (instead of the Japanese UI
# font) to mark them as "unset" - this introduces a new cellXfs
# style entry.
# ---------------------------------------------------------------

# Row 4: member 2 BAT -> BAT_2; member 3 slot cleared
$ws2.Range("B4").Value = "BAT_2"
$c = $ws2.Range("C4")
$c.Value = ""
$c.Font.Name = "Verdana"

# Row 5: member 4 slot cleared (D5)
$c = $ws2.Range("D5")
$c.Value = ""
$c.Font.Name = "Verdana"

# Row 6: member 2 BAT_2 -> WOLF; member 3 slot cleared
$ws2.Range("B6").Value = "WOLF"
$c = $ws2.Range("C6")
$c.Value = ""
$c.Font.Name = "Verdana"

# Row 7: every member slot cleared
foreach ($col in "B", "C", "D", "E", "F") {
    $c = $ws2.Range("$col`7")
    $c.Value = ""
    $c.Font.Name = "Verdana"
}

# Row 8: member 1 slot cleared
$c = $ws2.Range("B8")
$c.Value = ""
$c.Font.Name = "Verdana"
